$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.039.60'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '3.786.18'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").Value = '3.785.66'
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("E13").Value = '  -3.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '4.419.26'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").Value = '3.783.56'
$ws.Range("E16").Value = '  -2.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").Value = '67.921.24'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("E21").Value = '  -5.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '468.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("E24").Value = '  -8.74%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").Value = '3.933.25'
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").Value = '3.742.56'
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("E38").Value = '  -8.00%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.88%  '
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("E43").Value = '  -1.10%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '402.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.80%  '
$ws.Range("E49").Value = '  -7.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '40.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '142.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
